# Fix #300: the page count is incorrect when 'newPage'.asPagination() is set.
#
# The paragraph containing Page1/Page2/Page3/Page4 originally accumulates
# the page breaks inside the run that also carries the following page's
# text (e.g. one run holds two "<w:br w:type="page"/>" before "Page3").
# That duplication of page breaks is the bug: each run must carry exactly
# one page break, isolated in its own run (with an empty text node), and
# the following text must live in a run of its own.

$d = $word.ActiveDocument

# Locate the paragraph that contains the "Page1" marker text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Page1*Page2*Page3*Page4*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $full = $target.Range
    # Exclude the trailing paragraph mark so the <w:p> element (and its
    # rsid attributes) stay untouched; only the runs are replaced.
    $contentRange = $d.Range($full.Start, $full.End - 1)

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00B555E2"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Page1</w:t></w:r><w:r w:rsidRPr="00B555E2"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t/><w:br w:type="page"/></w:r><w:r w:rsidRPr="00B555E2"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Page2</w:t></w:r><w:r w:rsidRPr="00B555E2"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t/><w:br w:type="page"/></w:r><w:r w:rsidRPr="00B555E2"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Page3</w:t></w:r><w:r w:rsidRPr="00B555E2"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t/><w:br w:type="page"/></w:r><w:r w:rsidRPr="00B555E2"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Page4</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $contentRange.InsertXML($xml)
}
